$wb = $excel.ActiveWorkbook

# --- Sheet1: add the "{{ df2 | header }}" frame above the noheader/maxrows frame ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(6).Insert()
$ws1.Cells.Item(6, 1).Value = "{{ df2 | header }}"

# --- Sheet2 ("expected"): insert the rendered header row that goes with the new frame ---
$ws2 = $wb.Worksheets.Item("expected")
$ws2.Rows.Item(7).Insert()
$hdr = $ws2.Range("A3:E3").Value()
$ws2.Cells.Item(7, 1).Value = $hdr[1, 1]
$ws2.Cells.Item(7, 2).Value = $hdr[1, 2]
$ws2.Cells.Item(7, 3).Value = $hdr[1, 3]
$ws2.Cells.Item(7, 4).Value = $hdr[1, 4]
$ws2.Cells.Item(7, 5).Value = $hdr[1, 5]

# --- Selection / active sheet bookkeeping ---
$ws2.Range("B17").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A15").Select() | Out-Null
